# Auto update stock data
# Updates the "as of" date (column A) from 2025/12/31 to 2026/01/01 and the
# refreshed EBITDA figure (column B) for the first reporting row of each
# company block in the Stock Risk Scores sheet.
#
# Values are written as literal text (matching the workbook's existing
# inlineStr cells) rather than letting Excel auto-parse them as dates /
# numbers, by forcing a Text number format before the assignment and then
# restoring the cell to the "Normal" style so no stray formatting is left
# behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# row -> (new date, old EBITDA cell value, new EBITDA value)
$updates = @(
    @{ Row = 2;  EbitdaNew = "6.61" },
    @{ Row = 8;  EbitdaNew = "8.53" },
    @{ Row = 14; EbitdaNew = $null },
    @{ Row = 20; EbitdaNew = "12.69" },
    @{ Row = 26; EbitdaNew = "11.13" },
    @{ Row = 32; EbitdaNew = "27.27" },
    @{ Row = 38; EbitdaNew = $null },
    @{ Row = 44; EbitdaNew = "10.96" },
    @{ Row = 50; EbitdaNew = "11.07" },
    @{ Row = 56; EbitdaNew = "31.03" },
    @{ Row = 62; EbitdaNew = "11.39" },
    @{ Row = 68; EbitdaNew = "12.62" },
    @{ Row = 74; EbitdaNew = "16.43" }
)

foreach ($u in $updates) {
    $row = $u.Row
    Set-TextValue "A$row" "2026/01/01"
    if ($null -ne $u.EbitdaNew) {
        Set-TextValue "B$row" $u.EbitdaNew
    }
}
